# The source document currently contains a single paragraph (the
# "Abstract" body text). This script prepends the following block of
# new paragraphs in front of it:
#
#   Title:
#   Enhancing Psychological Assessment: PICA Survey & Mobile Self-Monitoring System
#   <empty paragraph>
#   Team Name:
#   PICA Software Development Unit
#   <empty paragraph>
#   Client Information:
#   Washington State University Psychology Clinic-Advancing personalized mental health assessment through innovative tools and technology.
#   <empty paragraph>
#   Abstract:
#
# We build each new paragraph from a literal WordprocessingML fragment
# via Range.InsertXML so the resulting markup (including an exactly
# empty <w:p/> for the blank separator paragraphs, and the preserved
# trailing space / "&" escaping on the title line) matches verbatim.

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-TextParagraphXml([string]$text, [bool]$preserveSpace) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    if ($preserveSpace) {
        return "<w:p xmlns:w='$wNs'><w:r><w:t xml:space='preserve'>$escaped</w:t></w:r></w:p>"
    }
    return "<w:p xmlns:w='$wNs'><w:r><w:t>$escaped</w:t></w:r></w:p>"
}

$emptyParagraphXml = "<w:p xmlns:w='$wNs'/>"

# Insert, from the bottom of the new block upward, always immediately
# before what is currently the document's first paragraph -- so the
# final top-to-bottom reading order comes out correct.

function Insert-ParagraphBeforeFirst([string]$xml) {
    $d = $word.ActiveDocument
    $target = $d.Paragraphs(1).Range
    $target.InsertParagraphBefore() | Out-Null
    $d.Paragraphs(1).Range.InsertXML($xml) | Out-Null
}

$emDash = [char]0x2014
$clientInfoText = "Washington State University Psychology Clinic" + $emDash + "Advancing personalized mental health assessment through innovative tools and technology."

Insert-ParagraphBeforeFirst (New-TextParagraphXml "Abstract:" $false)
Insert-ParagraphBeforeFirst $emptyParagraphXml
Insert-ParagraphBeforeFirst (New-TextParagraphXml $clientInfoText $false)
Insert-ParagraphBeforeFirst (New-TextParagraphXml "Client Information:" $false)
Insert-ParagraphBeforeFirst $emptyParagraphXml
Insert-ParagraphBeforeFirst (New-TextParagraphXml "PICA Software Development Unit" $false)
Insert-ParagraphBeforeFirst (New-TextParagraphXml "Team Name:" $false)
Insert-ParagraphBeforeFirst $emptyParagraphXml
Insert-ParagraphBeforeFirst (New-TextParagraphXml "Enhancing Psychological Assessment: PICA Survey & Mobile Self-Monitoring System " $true)
Insert-ParagraphBeforeFirst (New-TextParagraphXml "Title:" $false)
